$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$notes = $s.NotesPage
$notes.Shapes.Item(2).TextFrame.TextRange.Text = "Speaker notes here"
Write-Host "done"
